$p = $ppt.ActivePresentation
Write-Host "FullName=" $p.FullName
Write-Host "Path=" $p.Path
Write-Host "Name=" $p.Name
